$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.727484333333333
$ws.Range("H2").Value = 8.182453000000001
$ws.Range("I2").Value = 0.03096049453772388
$ws.Range("J2").Value = 0.03096049453772388
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.731629
$ws.Range("N2").Value = 8.194887
$ws.Range("O2").Value = 0.5547800938501829
$ws.Range("P2").Value = 0.554780093850183
$ws.Range("Q2").Value = 7.450475301979
$ws.Range("R2").Value = 67.054277717811
$ws.Range("S2").Value = 0.01717626606528653
$ws.Range("T2").Value = 0.01717626606528654
$ws.Range("G3").Value = 2.727484333333333
$ws.Range("H3").Value = 8.182453000000001
$ws.Range("I3").Value = 0.03096049453772388
$ws.Range("J3").Value = 0.03096049453772388
$ws.Range("M3").Value = 0.06813733333333333
$ws.Range("O3").Value = 0.01383834927121065
$ws.Range("P3").Value = 0.01383834927121065
$ws.Range("Q3").Value = 0.1858435091817778
$ws.Range("R3").Value = 1.672591582636
$ws.Range("S3").Value = 0.0004284421370224324
$ws.Range("T3").Value = 0.0004284421370224325
$ws.Range("G4").Value = 2.727484333333333
$ws.Range("H4").Value = 8.182453000000001
$ws.Range("I4").Value = 0.03096049453772388
$ws.Range("J4").Value = 0.03096049453772388
$ws.Range("M4").Value = 2.124038666666666
$ws.Range("N4").Value = 6.372115999999999
$ws.Range("O4").Value = 0.4313815568786064
$ws.Range("P4").Value = 0.4313815568786064
$ws.Range("Q4").Value = 5.793282186727555
$ws.Range("R4").Value = 52.139539680548
$ws.Range("S4").Value = 0.01335578633541492
$ws.Range("T4").Value = 0.01335578633541492
$ws.Range("I5").Value = 0.5986009007423507
$ws.Range("J5").Value = 0.5986009007423507
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.731629
$ws.Range("N5").Value = 8.194887
$ws.Range("O5").Value = 0.5547800938501829
$ws.Range("P5").Value = 0.554780093850183
$ws.Range("Q5").Value = 144.050064229082
$ws.Range("R5").Value = 1296.450578061738
$ws.Range("S5").Value = 0.3320918638926453
$ws.Range("T5").Value = 0.3320918638926454
$ws.Range("I6").Value = 0.5986009007423507
$ws.Range("J6").Value = 0.5986009007423507
$ws.Range("M6").Value = 0.06813733333333333
$ws.Range("O6").Value = 0.01383834927121065
$ws.Range("P6").Value = 0.01383834927121065
$ws.Range("Q6").Value = 3.593162630454222
$ws.Range("S6").Value = 0.008283648338533945
$ws.Range("T6").Value = 0.008283648338533947
$ws.Range("I7").Value = 0.5986009007423507
$ws.Range("J7").Value = 0.5986009007423507
$ws.Range("M7").Value = 2.124038666666666
$ws.Range("N7").Value = 6.372115999999999
$ws.Range("O7").Value = 0.4313815568786064
$ws.Range("P7").Value = 0.4313815568786064
$ws.Range("Q7").Value = 112.0093198448204
$ws.Range("R7").Value = 1008.083878603384
$ws.Range("S7").Value = 0.2582253885111714
$ws.Range("T7").Value = 0.2582253885111714
$ws.Range("G8").Value = 32.63402300000001
$ws.Range("H8").Value = 97.90206900000001
$ws.Range("I8").Value = 0.3704386047199253
$ws.Range("J8").Value = 0.3704386047199253
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.731629
$ws.Range("N8").Value = 8.194887
$ws.Range("O8").Value = 0.5547800938501829
$ws.Range("P8").Value = 0.554780093850183
$ws.Range("Q8").Value = 89.14404361346702
$ws.Range("R8").Value = 802.2963925212031
$ws.Range("S8").Value = 0.205511963892251
$ws.Range("T8").Value = 0.205511963892251
$ws.Range("G9").Value = 32.63402300000001
$ws.Range("H9").Value = 97.90206900000001
$ws.Range("I9").Value = 0.3704386047199253
$ws.Range("J9").Value = 0.3704386047199253
$ws.Range("M9").Value = 0.06813733333333333
$ws.Range("O9").Value = 0.01383834927121065
$ws.Range("P9").Value = 0.01383834927121065
$ws.Range("Q9").Value = 2.223595303158667
$ws.Range("R9").Value = 20.012357728428
$ws.Range("S9").Value = 0.005126258795654266
$ws.Range("T9").Value = 0.005126258795654267
$ws.Range("G10").Value = 32.63402300000001
$ws.Range("H10").Value = 97.90206900000001
$ws.Range("I10").Value = 0.3704386047199253
$ws.Range("J10").Value = 0.3704386047199253
$ws.Range("M10").Value = 2.124038666666666
$ws.Range("N10").Value = 6.372115999999999
$ws.Range("O10").Value = 0.4313815568786064
$ws.Range("P10").Value = 0.4313815568786064
$ws.Range("Q10").Value = 69.31592670088934
$ws.Range("R10").Value = 623.843340308004
$ws.Range("S10").Value = 0.15980038203202
$ws.Range("T10").Value = 0.15980038203202
